$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.920.15"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.239.53"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "269.79"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +3.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.49"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +15.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.635"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +7.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.56"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +6.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0957"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +4.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.40"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +20.73%  "
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.29"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +7.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.575.81"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +2.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.820"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +6.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.230.40"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.889.07"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.19"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +4.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.86"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.35"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -4.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.69"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.15"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +3.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  +7.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.50"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +12.16%  "
$ws.Range("E28").Value = "  +6.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.42"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -3.95%  "
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0914"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +5.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.01"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +2.86%  "
$ws.Range("E34").Value = "  +3.20%  "
$ws.Range("E35").Value = "  +2.27%  "
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0354"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.30"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -3.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.59"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +24.89%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.83"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -2.06%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.225"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +12.86%  "
$ws.Range("E42").Value = "  +3.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.52"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.43"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0998"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +1.64%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.43"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +1.08%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.41"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("E48").Value = "  +4.50%  "
$ws.Range("E49").Value = "  +2.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.454"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +3.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.458.16"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +2.08%  "
